$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 31499
$ws.Range("J3").Value = 31499
$ws.Range("L3").Value = 31499
$ws.Range("N3").Value = -31727
$ws.Range("H17").Value = 1914.3182
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2341
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 7023
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -7359
$ws.Range("H43").Value = 2541
$ws.Range("J43").Value = 2211.5
$ws.Range("L43").Value = 2211.5
$ws.Range("N43").Value = -2349.5
$ws.Range("H53").Value = 1170.7858
$ws.Range("I53").Value = 606.8
$ws.Range("J53").Value = 2580.75
$ws.Range("K53").Value = 606.8
$ws.Range("L53").Value = 2580.75
$ws.Range("M53").Value = 30.20000000000005
$ws.Range("N53").Value = -3854.75
$ws.Range("H70").Value = 5456.5
$ws.Range("I70").Value = 1826.5
$ws.Range("K70").Value = 5479.5
$ws.Range("M70").Value = -5209.5
$ws.Range("H73").Value = 5456.5
$ws.Range("I73").Value = 1826.5
$ws.Range("K73").Value = 5479.5
$ws.Range("M73").Value = -4543.5
$ws.Range("H102").Value = 31499
$ws.Range("J102").Value = 31499
$ws.Range("L102").Value = 31499
$ws.Range("N102").Value = -37989
$ws.Range("H105").Value = 35167.75
$ws.Range("J105").Value = 35167.75
$ws.Range("L105").Value = 35167.75
$ws.Range("N105").Value = -42155.75
$ws.Range("H106").Value = 7649.6665
$ws.Range("I106").Value = 7475
$ws.Range("J106").Value = 7999
$ws.Range("K106").Value = 7475
$ws.Range("L106").Value = 7999
$ws.Range("M106").Value = -6844
$ws.Range("N106").Value = -9261
$ws.Range("H112").Value = 2972.077
$ws.Range("J112").Value = 3053.0833
$ws.Range("L112").Value = 9159.249899999999
$ws.Range("N112").Value = -11375.2499
$ws.Range("H113").Value = 1668.3334
$ws.Range("I113").Value = 1668.3334
$ws.Range("K113").Value = 1668.3334
$ws.Range("M113").Value = 1585.6666
$ws.Range("H138").Value = 2234
$ws.Range("J138").Value = 2205.5
$ws.Range("L138").Value = 6616.5
$ws.Range("N138").Value = -16896.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6471.3335
$ws.Range("I74").Value = 3546.5
$ws.Range("K74").Value = 3546.5
$ws.Range("M74").Value = -2672.5
$ws.Range("H77").Value = 6471.3335
$ws.Range("I77").Value = 3546.5
$ws.Range("K77").Value = 17732.5
$ws.Range("M77").Value = -13364.5
$ws.Range("H122").Value = 2037
$ws.Range("I122").Value = 2037
$ws.Range("K122").Value = 6111
$ws.Range("M122").Value = -3661
$ws.Range("H132").Value = 5429.5
$ws.Range("I132").Value = 804.5
$ws.Range("K132").Value = 2413.5
$ws.Range("M132").Value = 116.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3163.8333
$ws.Range("I134").Value = 2730
$ws.Range("K134").Value = 8190
$ws.Range("M134").Value = -5655

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1425.6666
$ws.Range("J16").Value = 1877
$ws.Range("L16").Value = 1877
$ws.Range("N16").Value = -2451
$ws.Range("H22").Value = 3361.75
$ws.Range("I22").Value = 2723.5
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 2723.5
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -2373.5
$ws.Range("N22").Value = -4700
$ws.Range("H31").Value = 7397.7827
$ws.Range("I31").Value = 3637.1667
$ws.Range("J31").Value = 8725.058999999999
$ws.Range("K31").Value = 3637.1667
$ws.Range("L31").Value = 8725.058999999999
$ws.Range("M31").Value = -3342.1667
$ws.Range("N31").Value = -9315.058999999999
$ws.Range("H34").Value = 7397.7827
$ws.Range("I34").Value = 3637.1667
$ws.Range("J34").Value = 8725.058999999999
$ws.Range("K34").Value = 3637.1667
$ws.Range("L34").Value = 8725.058999999999
$ws.Range("M34").Value = -3435.1667
$ws.Range("N34").Value = -9129.058999999999
$ws.Range("H38").Value = 15234.25
$ws.Range("J38").Value = 15234.25
$ws.Range("L38").Value = 15234.25
$ws.Range("N38").Value = -15988.25
$ws.Range("H46").Value = 15234.25
$ws.Range("J46").Value = 15234.25
$ws.Range("L46").Value = 15234.25
$ws.Range("N46").Value = -15656.25
$ws.Range("H88").Value = 11513.667
$ws.Range("J88").Value = 14000
$ws.Range("L88").Value = 14000
$ws.Range("N88").Value = -14812
$ws.Range("H91").Value = 11513.667
$ws.Range("J91").Value = 14000
$ws.Range("L91").Value = 14000
$ws.Range("N91").Value = -16808
$ws.Range("H106").Value = 35975
$ws.Range("J106").Value = 35975
$ws.Range("L106").Value = 35975
$ws.Range("N106").Value = -38499
$ws.Range("H107").Value = 355.7857
$ws.Range("I107").Value = 583.3333
$ws.Range("K107").Value = 583.3333
$ws.Range("M107").Value = 1336.6667
$ws.Range("H113").Value = 1425.6666
$ws.Range("J113").Value = 1877
$ws.Range("L113").Value = 1877
$ws.Range("N113").Value = -6217
$ws.Range("H134").Value = 4221
$ws.Range("I134").Value = 3534.1428
$ws.Range("J134").Value = 6625
$ws.Range("K134").Value = 10602.4284
$ws.Range("L134").Value = 19875
$ws.Range("M134").Value = -8067.428400000001
$ws.Range("N134").Value = -24945

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 100.5
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H131").Value = 1992.5625
$ws.Range("I131").Value = 1648.2
$ws.Range("J131").Value = 2566.5
$ws.Range("K131").Value = 4944.6
$ws.Range("L131").Value = 7699.5
$ws.Range("M131").Value = 95.39999999999964
$ws.Range("N131").Value = -17779.5
$ws.Range("H134").Value = 6483.1113
$ws.Range("I134").Value = 1225
$ws.Range("J134").Value = 16999.334
$ws.Range("K134").Value = 3675
$ws.Range("L134").Value = 50998.00199999999
$ws.Range("M134").Value = 1395
$ws.Range("N134").Value = -61138.00199999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 35599.8
$ws.Range("J15").Value = 35599.8
$ws.Range("L15").Value = 35599.8
$ws.Range("N15").Value = -36175.8
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 35599.8
$ws.Range("J81").Value = 35599.8
$ws.Range("L81").Value = 35599.8
$ws.Range("N81").Value = -37595.8
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 35599.8
$ws.Range("J84").Value = 35599.8
$ws.Range("L84").Value = 106799.4
$ws.Range("N84").Value = -116783.4
$ws.Range("H122").Value = 387624.7
$ws.Range("I122").Value = 419278.66
$ws.Range("K122").Value = 1257835.98
$ws.Range("M122").Value = -1255385.98
$ws.Range("H132").Value = 67625.75
$ws.Range("I132").Value = 130352.375
$ws.Range("J132").Value = 4899.125
$ws.Range("K132").Value = 391057.125
$ws.Range("L132").Value = 14697.375
$ws.Range("M132").Value = -388527.125
$ws.Range("N132").Value = -19757.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 7087.375
$ws.Range("I68").Value = 2233.3333
$ws.Range("J68").Value = 9999.799999999999
$ws.Range("K68").Value = 2233.3333
$ws.Range("L68").Value = 9999.799999999999
$ws.Range("M68").Value = -1484.3333
$ws.Range("N68").Value = -11497.8
$ws.Range("H71").Value = 7087.375
$ws.Range("I71").Value = 2233.3333
$ws.Range("J71").Value = 9999.799999999999
$ws.Range("K71").Value = 11166.6665
$ws.Range("L71").Value = 49999
$ws.Range("M71").Value = -7422.666499999999
$ws.Range("N71").Value = -57487
$ws.Range("H100").Value = 8020.5415
$ws.Range("I100").Value = 7943.2
$ws.Range("J100").Value = 8040.8945
$ws.Range("K100").Value = 7943.2
$ws.Range("L100").Value = 8040.8945
$ws.Range("M100").Value = -7402.2
$ws.Range("N100").Value = -9122.8945
$ws.Range("H122").Value = 7597.4
$ws.Range("I122").Value = 6997
$ws.Range("K122").Value = 20991
$ws.Range("M122").Value = -18541
$ws.Range("H136").Value = 1372.25
$ws.Range("I136").Value = 1161.3334
$ws.Range("K136").Value = 3484.0002
$ws.Range("M136").Value = -934.0001999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 15000
$ws.Range("J31").Value = 15000
$ws.Range("L31").Value = 15000
$ws.Range("N31").Value = -15696
$ws.Range("H107").Value = 585.2857
$ws.Range("I107").Value = 499.45456
$ws.Range("K107").Value = 1498.36368
$ws.Range("M107").Value = 421.6363200000001
$ws.Range("H122").Value = 2630.7058
$ws.Range("J122").Value = 3986.2222
$ws.Range("L122").Value = 11958.6666
$ws.Range("N122").Value = -16858.6666
$ws.Range("H132").Value = 4059.4
$ws.Range("I132").Value = 2515.6667
$ws.Range("K132").Value = 7547.000100000001
$ws.Range("M132").Value = -5017.000100000001
$ws.Range("H136").Value = 3742.8845
$ws.Range("I136").Value = 2324.6
$ws.Range("J136").Value = 5676.909
$ws.Range("K136").Value = 6973.799999999999
$ws.Range("L136").Value = 17030.727
$ws.Range("M136").Value = -4423.799999999999
$ws.Range("N136").Value = -22130.727
